# Apply updated cryptos list values (prices, volumes, and row reorders)
# Source: GitHub Actions scheduled update, Thu Jan 25 08:54:01 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "40.129.77"
$ws.Range("E2").Value = "  +0.43%  "

# Row 3
$ws.Range("D3").Value = "2.224.12"
$ws.Range("E3").Value = "  -0.51%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "291.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.55%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.55%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  +0.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.26%  "

# Row 12
$ws.Range("E12").Value = "  +3.28%  "

# Row 13
$ws.Range("E13").Value = "  +1.16%  "

# Row 14
$ws.Range("D14").Value = "2.568.27"
$ws.Range("E14").Value = "  -0.37%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.74%  "

# Row 16
$ws.Range("D16").Value = "2.224.29"
$ws.Range("E16").Value = "  -0.06%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.728"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
$ws.Range("D18").Value = "40.082.48"
$ws.Range("E18").Value = "  +0.53%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0886"
$ws.Range("E19").Value = "  -0.92%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.13%  "

# Row 21
$ws.Range("E21").Value = "  +0.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("E24").Value = "  -0.05%  "

# Row 25
$ws.Range("E25").Value = "  +1.66%  "

# Row 26
$ws.Range("E26").Value = "  -0.42%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.16%  "

# Row 28
$ws.Range("E28").Value = "  -1.17%  "

# Row 29
$ws.Range("E29").Value = "  -0.68%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "156.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.27%  "

# Row 32
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("E33").Value = "  +1.75%  "

# Row 34
$ws.Range("E34").Value = "  +1.16%  "

# Row 35
$ws.Range("E35").Value = "  -1.22%  "

# Row 36
$ws.Range("E36").Value = "  +6.48%  "

# Row 37
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("E38").Value = "  -5.34%  "

# Row 39
$ws.Range("E39").Value = "  -2.02%  "

# Row 40
$ws.Range("E40").Value = "  +1.61%  "

# Row 41
$ws.Range("D41").Value = "2.116.61"
$ws.Range("E41").Value = "  +7.78%  "

# Row 42
$ws.Range("E42").Value = "  +1.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.00%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0268"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.17%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.32%  "

# Row 47
$ws.Range("E47").Value = "  +2.79%  "

# Row 48
$ws.Range("D48").Value = "2.435.45"
$ws.Range("E48").Value = "  -0.55%  "

# Row 49
$ws.Range("E49").Value = "  -0.40%  "

# Row 50
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.96%  "

# Row 51
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "69.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.97%  "
